$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Update the date in A1 (45310 -> 45311)
$ws.Range("A1").Value = 45311

# Update prices in D29 and D30
$ws.Range("D29").Value = 264.2
$ws.Range("D30").Value = 179.23
